# CCDC_Datasets.xlsx edit script
# Implements: new "ResourceCode" column (dbGaP / TARGET) inserted after DatasetName,
# and a new "grantInfo" column holding the full grant text (moved out of the
# "grant" column, with a couple of trailing spaces trimmed) on the
# CoreAdditional sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("CoreAdditional")

# ---------------------------------------------------------------------------
# 1. Insert a new column C ("ResourceCode") - this shifts the existing
#    numOfCases..grant columns (C..J) one slot to the right (D..K), carrying
#    their values/styles with them.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(3).Insert()

$ws2.Range("C1").Value = "ResourceCode"
$ws2.Range("C2").Value = "dbGaP"
$ws2.Range("C3").Value = "TARGET"

# ---------------------------------------------------------------------------
# 2. The caseAge cell in row 2 (now column F) used to carry a one-off
#    "Lato 8pt" font style (the font is being retired); restore it to the
#    plain wrap-text formatting used by the same column in row 3.
# ---------------------------------------------------------------------------
$ws2.Range("F2").Style = "Normal"
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Add a new "grantInfo" column (L) holding the full grant text that used
#    to live in the "grant" column (K, after the shift above) on row 3; trim
#    the stray trailing spaces in that text and clear it out of column K.
# ---------------------------------------------------------------------------
$ws2.Range("L1").Value = "grantInfo"
$ws2.Range("L1").NumberFormat = "@"

$grantText = "261200800001E-12-0-40`nTherapeutically Applicable Research to Generate Effective Treatments (TARGET)`nHHSN261200800001E`nNCI Contract`nU10CA180886`nCOG NCTN Network Group Operations Center"
$ws2.Range("L3").Value = $grantText
$ws2.Range("L3").NumberFormat = "@"
$ws2.Range("L3").WrapText = $true

$ws2.Range("K3").ClearContents()

# ---------------------------------------------------------------------------
# 4. Row-height touch-up to reflect the extra wrapped content.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(2).RowHeight = 30
$ws2.Rows.Item(3).RowHeight = 150
